$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI ligand-receptor pair table: "M2" and "sCs" sending clusters
# each paired against target clusters ECs, FAPs, M2, sCs (per Dr Hou advice).
$data = @(
    ,("M2", "Cd22", "Ptprc", "ECs", 3, 1, 14.445741, 43.337223, 0.9879457199603994, 0.9879457199603994, 2, 0.6666666666666666, 0.1546356666666667, 0.463907, 0.0006541814359458435, 0.0006541814359458435, 2.233826790029, 20.104441110261, 0.0006462957497202443, 0.0006462957497202443)
    ,("M2", "Cd22", "Ptprc", "FAPs", 3, 1, 14.445741, 43.337223, 0.9879457199603994, 0.9879457199603994, 3, 1, 0.3253436666666666, 0.976031, 0.001376356384162467, 0.001376356384162467, 4.699830344656999, 42.298473101913, 0.00135976539887348, 0.00135976539887348)
    ,("M2", "Cd22", "Ptprc", "M2", 3, 1, 14.445741, 43.337223, 0.9879457199603994, 0.9879457199603994, 3, 1, 231.9004236666667, 695.701271, 0.9810476161216116, 0.9810476161216116, 3349.973458078937, 30149.76112271044, 0.9692217934246991, 0.9692217934246991)
    ,("M2", "Cd22", "Ptprc", "sCs", 3, 1, 14.445741, 43.337223, 0.9879457199603994, 0.9879457199603994, 3, 1, 3.999992666666667, 11.999978, 0.01692184605828007, 0.01692184605828007, 57.78285806456601, 520.0457225810941, 0.01671786538710655, 0.01671786538710655)
    ,("sCs", "Cd22", "Ptprc", "ECs", 1, 0.3333333333333333, 0.1762576666666667, 0.528773, 0.01205428003960061, 0.01205428003960061, 2, 0.6666666666666666, 0.1546356666666667, 0.463907, 0.0006541814359458435, 0.0006541814359458435, 0.02725572179011111, 0.245301496111, 0.000007885686225599245, 0.000007885686225599245)
    ,("sCs", "Cd22", "Ptprc", "FAPs", 1, 0.3333333333333333, 0.1762576666666667, 0.528773, 0.01205428003960061, 0.01205428003960061, 3, 1, 0.3253436666666666, 0.976031, 0.001376356384162467, 0.001376356384162467, 0.05734431555144445, 0.516098839963, 0.00001659098528898649, 0.00001659098528898649)
    ,("sCs", "Cd22", "Ptprc", "M2", 1, 0.3333333333333333, 0.1762576666666667, 0.528773, 0.01205428003960061, 0.01205428003960061, 3, 1, 231.9004236666667, 695.701271, 0.9810476161216116, 0.9810476161216116, 40.87422757449811, 367.8680481704831, 0.0118258226969125, 0.0118258226969125)
    ,("sCs", "Cd22", "Ptprc", "sCs", 1, 0.3333333333333333, 0.1762576666666667, 0.528773, 0.01205428003960061, 0.01205428003960061, 3, 1, 3.999992666666667, 11.999978, 0.01692184605828007, 0.01692184605828007, 0.7050293741104445, 6.345264366994001, 0.0002039806711735196, 0.0002039806711735196)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
